$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.612.69"
$ws.Range("E2").Value = "  -0.42%  "

$ws.Range("D3").Value = "3.443.84"
$ws.Range("E3").Value = "  -1.13%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "'408.23"
$ws.Range("E5").Value = "  -0.41%  "

$ws.Range("E6").Value = "  +0.85%  "

$ws.Range("E7").Value = "  -1.60%  "

$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.00%  "

$ws.Range("D9").Value = "'0.689"
$ws.Range("E9").Value = "  -1.10%  "

$ws.Range("D10").Value = "'0.124"
$ws.Range("E10").Value = "  -3.92%  "

$ws.Range("D11").Value = "'42.22"
$ws.Range("E11").Value = "  -2.46%  "

$ws.Range("E12").Value = "  -0.81%  "

$ws.Range("D13").Value = "'8.48"
$ws.Range("E13").Value = "  -3.11%  "

$ws.Range("E14").Value = "  -1.50%  "

$ws.Range("D15").Value = "3.523.81"
$ws.Range("E15").Value = "  +1.38%  "

$ws.Range("D16").Value = "62.523.94"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").Value = "'11.42"
$ws.Range("E17").Value = "  +5.02%  "

$ws.Range("E18").Value = "  -2.30%  "

$ws.Range("E19").Value = "  -4.68%  "

$ws.Range("E20").Value = "  -5.64%  "

$ws.Range("D21").Value = "'84.07"
$ws.Range("E21").Value = "  +1.41%  "

$ws.Range("D22").Value = "'314.26"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("E23").Value = "  -1.63%  "

$ws.Range("D24").Value = "'3.17"
$ws.Range("E24").Value = "  -0.36%  "

$ws.Range("E25").Value = "  +8.02%  "

$ws.Range("D26").Value = "'29.83"
$ws.Range("E26").Value = "  -2.23%  "

$ws.Range("D27").Value = "'8.24"
$ws.Range("E27").Value = "  -1.13%  "

$ws.Range("D28").Value = "'2.80"
$ws.Range("E28").Value = "  +5.60%  "

$ws.Range("D29").Value = "'7.59"
$ws.Range("E29").Value = "  -2.36%  "

$ws.Range("D30").Value = "'0.174"
$ws.Range("E30").Value = "  -3.45%  "

$ws.Range("E31").Value = "  -3.89%  "

$ws.Range("D32").Value = "'42.39"
$ws.Range("E32").Value = "  -1.67%  "

$ws.Range("E33").Value = "  -0.07%  "

$ws.Range("D34").Value = "'11.41"
$ws.Range("E34").Value = "  -4.48%  "

$ws.Range("E35").Value = "  -1.96%  "

$ws.Range("D36").Value = "'51.45"
$ws.Range("E36").Value = "  -1.86%  "

$ws.Range("D37").Value = "'0.997"
$ws.Range("E37").Value = "  -0.07%  "

$ws.Range("E38").Value = "  -5.93%  "

$ws.Range("D39").Value = "'0.326"
$ws.Range("E39").Value = "  +13.39%  "

$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("D41").Value = "'138.28"
$ws.Range("E41").Value = "  -0.10%  "

$ws.Range("E42").Value = "  -0.42%  "

$ws.Range("E43").Value = "  -0.65%  "

$ws.Range("E44").Value = "  +0.41%  "

$ws.Range("E45").Value = "  -4.46%  "

$ws.Range("D47").Value = "'21.31"
$ws.Range("E47").Value = "  -5.32%  "

$ws.Range("D48").Value = "2.127.08"
$ws.Range("E48").Value = "  -3.61%  "

$ws.Range("E49").Value = "  -3.40%  "

$ws.Range("B50").Value = "ThetaToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D50").Value = "'1.93"
$ws.Range("E50").Value = "  +2.60%  "

$ws.Range("B51").Value = "Fetch.AI"
$ws.Range("C51").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D51").Value = "'1.74"
$ws.Range("E51").Value = "  +20.81%  "
